# Append the latest metric reading as a new row at the bottom of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the next empty row right after the current data (row 49 -> row 50).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# New data point: timestamp (stored as text, like the existing rows) and metric value.
$ws.Cells.Item($newRow, 1).Value = "2025-04-29 05:49:57"
$ws.Cells.Item($newRow, 2).Value = 160
